# Update test data on the TESTDATA sheet: row 8 ("postProductWithoutPOJO")
# had placeholder name/type values "amuthan" / "sakthivel" - replace them
# with "Manjeet" / "Singh".
$wb = $excel.ActiveWorkbook
$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsTestData = $wb.Worksheets.Item("TESTDATA")

$wsTestData.Range("F8").Value = "Manjeet"
$wsTestData.Range("G8").Value = "Singh"

# Match the author's recorded selection/view state from the commit: the
# RUNMANAGER sheet was left with the cursor at A32 (no longer the active
# tab) while TESTDATA became the active tab, selected at G8 (where the
# edit was made).
$wsRunManager.Activate()
$wsRunManager.Range("A32").Select()

$wsTestData.Activate()
$wsTestData.Range("G8").Select()
